# Add a new "2021-01-11" attendance sheet, modeled on the existing
# "2021-01-07" sheet (same header row / column-A styling / page margins),
# populated with that day's readings.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("2021-01-07")

# Duplicate the source sheet (carries over sheetPr, styles, column-A /
# header formatting and page margins) and drop it right after the source
# sheet, i.e. as the new last tab.
$sourceSheet.Copy($null, $sourceSheet)
$ws = $wb.Worksheets.Item($sourceSheet.Index + 1)
$ws.Name = "2021-01-11"

# Attendance readings for 2021-01-11 (Sr.No, Name, Address, Job,
# Time-Stamp, SpO2_value, Heart-rate, Compensated, Ambient).
$data = @(
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "18:35:14", 97.59385534014351, 0, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:31:51", 97.66750902355957, 66.57590464616032, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:36:47", 97.26071028597477, 85.7017862963194, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:40:10", 97.38368981463954, 70.06637470606302, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:44:27", 96.9643448485741, 131.0042587766994, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:44:46", 97.66721486857014, 124.4869430863692, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:46:16", 97.2858232863867, 51.98913957578586, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:47:30", 97.47390529427025, 94.44407387681979, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "19:47:59", 97.60036380894805, 74.83467442504315, "NA", "NA")
)

# Row 1 (header) + row 2's formatting/values already came across with the
# sheet copy; overwrite row 2 with the real data and append rows 3-10.
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}

# Column A keeps the bold/bordered/centered style on every data row, same
# as the source sheet - replicate that formatting down the new rows.
$ws.Range("A2").Copy()
$ws.Range("A3:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Copying a sheet makes it active; restore Sheet1 as the selected tab so
# the workbook-level active sheet/tab stay as they were.
$wb.Worksheets.Item(1).Activate()
